# Applies the 4 textual corrections described by the commit:
#  1. Slide 1  - Report Date textbox: 2026-02-11 -> 2026-02-12
#  2. Slide 30 - Mobile Key Message: drop "EUR " before both amounts
#  3. Slide 31 - Fixed Key Message: drop "EUR " before the amount
#  4. Slide 32 - B2B Key Message: drop "EUR " before the amount

$p = $ppt.ActivePresentation

# --- Slide 1: update the Report Date value -------------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(15).TextFrame.TextRange.Text = "2026-02-12"

# --- Slide 30: Mobile look - Key Message textbox --------------------------
$s30 = $p.Slides.Item(30)
$s30.Shapes.Item(21).TextFrame.TextRange.Text = "Key Message: Mobile service revenue at 1520.0M; up 2.8% YoY; ARPU 12.8"

# --- Slide 31: Fixed look - Key Message textbox ----------------------------
$s31 = $p.Slides.Item(31)
$s31.Shapes.Item(21).TextFrame.TextRange.Text = "Key Message: Fixed service revenue 795.0M; growth -1.1% YoY; Fiber subs 1480K"

# --- Slide 32: B2B look - Key Message textbox ------------------------------
$s32 = $p.Slides.Item(32)
$s32.Shapes.Item(21).TextFrame.TextRange.Text = "Key Message: B2B revenue 520.0M; growth +8.5% YoY; 16.8% of total revenue"
